$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.386.06'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").Value = '3.378.56'
$ws.Range("E3").Value = '  -1.98%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'567.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.25%  '

$ws.Range("D6").Value = "'139.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.29%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.381.68'
$ws.Range("E8").Value = '  -1.92%  '

$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("E10").Value = '  -3.53%  '

$ws.Range("E11").Value = '  -2.95%  '

$ws.Range("E12").Value = '  -1.29%  '

$ws.Range("D13").Value = '3.953.81'
$ws.Range("E13").Value = '  -2.06%  '

$ws.Range("E14").Value = '  +0.93%  '

$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").Value = '3.378.77'
$ws.Range("E16").Value = '  -2.07%  '

$ws.Range("D17").Value = "'0.0000170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.34%  '

$ws.Range("D18").Value = '60.525.03'
$ws.Range("E18").Value = '  -1.99%  '

$ws.Range("D19").Value = "'6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.53%  '

$ws.Range("D20").Value = "'13.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.03%  '

$ws.Range("E21").Value = '  -5.39%  '

$ws.Range("D22").Value = "'386.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("D23").Value = "'0.553"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.28%  '

$ws.Range("E24").Value = '  -0.17%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("E26").Value = '  -8.47%  '

$ws.Range("D27").Value = '3.520.40'
$ws.Range("E27").Value = '  -1.91%  '

$ws.Range("E28").Value = '  -1.71%  '

$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").Value = "'7.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.69%  '

$ws.Range("D31").Value = "'7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.00%  '

$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("D33").Value = "'1.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.55%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = "'23.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.51%  '

$ws.Range("D36").Value = '3.408.13'
$ws.Range("E36").Value = '  -1.84%  '

$ws.Range("D38").Value = "'6.86"
$ws.Range("D38").Style = "Normal"

$ws.Range("E39").Value = '  -5.12%  '

$ws.Range("E40").Value = '  -4.92%  '

$ws.Range("D41").Value = "'0.0768"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.23%  '

$ws.Range("D42").Value = "'27.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D44").Value = "'0.779"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.65%  '

$ws.Range("D45").Value = "'4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.28%  '

$ws.Range("D46").Value = "'41.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.22%  '

$ws.Range("D48").Value = '2.513.94'
$ws.Range("E48").Value = '  -3.56%  '

$ws.Range("E49").Value = '  -4.51%  '

$ws.Range("D50").Value = "'23.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("E51").Value = '  -3.33%  '
